$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep a text value (matching the source file's
    # inlineStr cell type) even when the text looks like a number,
    # then restore the default "Normal" style so no stray formatting
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "29.190.51"
$ws.Range("E2").Value = "  +3.12%  "

# Row 3
$ws.Range("D3").Value = "1.581.23"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "212.37"
$ws.Range("E5").Value = "  +1.18%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.510"
$ws.Range("E6").Value = "  +5.59%  "

# Row 7
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
Set-TextValue $ws.Range("D8") "26.20"
$ws.Range("E8").Value = "  +10.19%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.249"
$ws.Range("E9").Value = "  +2.53%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0594"
$ws.Range("E10").Value = "  +1.96%  "

# Row 11
$ws.Range("E11").Value = "  +1.61%  "

# Row 12
$ws.Range("D12").Value = "1.807.58"
$ws.Range("E12").Value = "  +1.87%  "

# Row 13
$ws.Range("D13").Value = "1.578.11"
$ws.Range("E13").Value = "  +1.62%  "

# Row 14
$ws.Range("D14").Value = "29.185.13"
$ws.Range("E14").Value = "  +3.17%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.523"
$ws.Range("E15").Value = "  +2.68%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "3.71"
$ws.Range("E16").Value = "  +2.66%  "

# Row 17
Set-TextValue $ws.Range("D17") "62.24"

# Row 18
Set-TextValue $ws.Range("D18") "236.06"
$ws.Range("E18").Value = "  +3.86%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.44"
$ws.Range("E19").Value = "  +1.70%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  +2.85%  "

# Row 21
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("E22").Value = "  +2.09%  "

# Row 23
$ws.Range("E23").Value = "  +3.14%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.09"
$ws.Range("E24").Value = "  +3.26%  "

# Row 25
Set-TextValue $ws.Range("D25") "153.81"
$ws.Range("E25").Value = "  +1.50%  "

# Row 26
Set-TextValue $ws.Range("D26") "15.16"
$ws.Range("E26").Value = "  +2.83%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.108"
$ws.Range("E27").Value = "  +4.39%  "

# Row 28
$ws.Range("E28").Value = "  +1.69%  "

# Row 29
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0468"
$ws.Range("E30").Value = "  +0.28%  "

# Row 31
$ws.Range("E31").Value = "  +0.55%  "

# Row 32
$ws.Range("E32").Value = "  +1.68%  "

# Row 33
$ws.Range("D33").Value = "1.421.63"
$ws.Range("E33").Value = "  +2.43%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.08"
$ws.Range("E34").Value = "  +2.17%  "

# Row 35
$ws.Range("E35").Value = "  -1.79%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.51"
$ws.Range("E36").Value = "  +1.90%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.77"
$ws.Range("E37").Value = "  +7.15%  "

# Row 38
$ws.Range("E38").Value = "  -1.68%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0164"
$ws.Range("E39").Value = "  +1.62%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.530"
$ws.Range("E40").Value = "  +3.54%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.98"
$ws.Range("E41").Value = "  +2.56%  "

# Row 42
Set-TextValue $ws.Range("D42") "53.01"
$ws.Range("E42").Value = "  +25.24%  "

# Row 43
$ws.Range("E43").Value = "  -0.17%  "

# Row 44
$ws.Range("E44").Value = "  +1.60%  "

# Row 45
$ws.Range("E45").Value = "  +0.35%  "

# Row 46
Set-TextValue $ws.Range("D46") "64.59"
$ws.Range("E46").Value = "  +4.55%  "

# Row 47
$ws.Range("E47").Value = "  -0.37%  "

# Row 48
$ws.Range("D48").Value = "1.718.80"
$ws.Range("E48").Value = "  +1.82%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.844"
$ws.Range("E49").Value = "  -6.70%  "

# Row 50
Set-TextValue $ws.Range("D50") "85.39"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("E51").Value = "  -1.19%  "
